# Writing Progress Tracker - apply tracked changes
# - Fill in "Type" column for a few existing papers
# - Backfill "Date Reviewed" for rows 4 & 5
# - Add four new papers (rows 20-23) to the "Papers List" sheet
# - Resize columns to fit the new "Source" column that was inserted
# - Update the sheet's zoom / selection to match the final saved state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Papers List")

# --- Fill in the "Type" column (C) for a few already-existing rows ---
$ws.Range("C3").Value = "Review"
$ws.Range("C4").Value = "Review"

# --- Backfill "Date Reviewed" (H) for rows 4 & 5 to match "Date Found" (G) ---
$ws.Range("H4").NumberFormat = $ws.Range("G4").NumberFormat
$ws.Range("H4").Value = $ws.Range("G4").Value()

$ws.Range("H5").NumberFormat = $ws.Range("G5").NumberFormat
$ws.Range("H5").Value = $ws.Range("G5").Value()

# This "Lab"-type paper is entered in between the two new-paper rows below
# (matches the order the underlying shared-string table was built in)

# --- New paper row 20 ---
$ws.Range("A20").Value = "Mixtures of Chemical Pollutants at European Legislation Safety Concentrations: How Safe Are They? "
$ws.Range("B20").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B20").Value = 41883
$ws.Range("E20").Value = "Carvalho et al."
$ws.Range("D20").Value = "Toxicological Sciences"
$ws.Range("F20").Value = "Mendely"
$ws.Range("G20").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("G20").Value = 43255

# --- New paper row 21 ---
$ws.Range("E21").Value = "C. Ritz, J. Streiberg"
$ws.Range("A21").Value = "From additivity to synergism " + [char]0x2013 + " A modelling perspective"
$ws.Range("B21").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B21").Value = 41883
$ws.Range("C21").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("C21").Value = "Technique"
$ws.Range("D21").Value = "Synergy"
$ws.Range("F21").Value = "Mendely"
$ws.Range("G21").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("G21").Value = 43255

# --- Type for an older row (row 5), entered around the same time as the above ---
$ws.Range("C5").Value = "Lab"

# --- New paper row 22 ---
$ws.Range("A22").Value = "Mechanisms of nickel toxicity in microorganisms."
$ws.Range("B22").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B22").Value = 40848
$ws.Range("C22").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("C22").Value = "Review"
$ws.Range("E22").Value = "L. Macomber, R. Hausinger"
$ws.Range("D22").Value = "Metallomics"
$ws.Range("F22").Value = "Google"
$ws.Range("G22").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("G22").Value = 43256

# --- New paper row 23 (title only, not yet reviewed) ---
$ws.Range("A23").Value = "Review of the molluscicide metaldehyde in the environment"
$ws.Range("G23").NumberFormat = $ws.Range("G2").NumberFormat
$ws.Range("G23").Value = 43256

# --- Column widths: "Source" (F) is a new column squeezed out of the old
#     Journal/Authors columns, and "Date Reviewed" (H) shrinks to fit ---
$ws.Columns.Item(4).ColumnWidth = 47.0
$ws.Columns.Item(5).ColumnWidth = 33.333333333333336
$ws.Columns.Item(6).ColumnWidth = 45.166666666666664
$ws.Columns.Item(8).ColumnWidth = 13.666666666666666

# --- View state: zoom to 100% and select the last entry cell ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("G23").Select()
